# Weekly fruit/vegetable price update.
# Three new price records (rows) are inserted at the top of the data block
# (right after the header's first existing entries), shifting all the
# subsequent rows down by three. The three new rows contain fresh
# "Camote" (Zapallo) price observations for "Región Metropolitana".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 16; this pushes old rows 16-62 down to 19-65
$ws.Rows("16:18").Insert()

# Copy the formatting/style of the row directly below (the shifted former
# row 16, now at row 19) into the three new rows so number formats
# (e.g. the date style on column D) match the rest of the table.
$ws.Rows("19").Copy()
$ws.Rows("16").PasteSpecial()
$ws.Rows("19").Copy()
$ws.Rows("17").PasteSpecial()
$ws.Rows("19").Copy()
$ws.Rows("18").PasteSpecial()

# --- New row 16 ---
$ws.Range("A16").Value2 = 1
$ws.Range("B16").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C16").Value2 = "Arica y Parinacota"
$ws.Range("D16").Value2 = 45054
$ws.Range("E16").Value2 = 15
$ws.Range("F16").Value2 = 100112045
$ws.Range("G16").Value2 = "Zapallo"
$ws.Range("H16").Value2 = "Camote"
$ws.Range("I16").Value2 = "1a (cosecha)"
$ws.Range("J16").Value2 = 850
$ws.Range("K16").Value2 = 340
$ws.Range("L16").Value2 = 350
$ws.Range("M16").Value2 = 346
$ws.Range("N16").Value2 = "`$/kilo (volumen en unidades)"
$ws.Range("O16").Value2 = "Región Metropolitana"
$ws.Range("P16").Value2 = 346
$ws.Range("Q16").Value2 = 1
$ws.Range("R16").Value2 = "Hortaliza"

# --- New row 17 ---
$ws.Range("A17").Value2 = 1
$ws.Range("B17").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C17").Value2 = "Arica y Parinacota"
$ws.Range("D17").Value2 = 45054
$ws.Range("E17").Value2 = 15
$ws.Range("F17").Value2 = 100112045
$ws.Range("G17").Value2 = "Zapallo"
$ws.Range("H17").Value2 = "Camote"
$ws.Range("I17").Value2 = "2a (cosecha)"
$ws.Range("J17").Value2 = 700
$ws.Range("K17").Value2 = 320
$ws.Range("L17").Value2 = 330
$ws.Range("M17").Value2 = 326
$ws.Range("N17").Value2 = "`$/kilo (volumen en unidades)"
$ws.Range("O17").Value2 = "Región Metropolitana"
$ws.Range("P17").Value2 = 326
$ws.Range("Q17").Value2 = 1
$ws.Range("R17").Value2 = "Hortaliza"

# --- New row 18 ---
$ws.Range("A18").Value2 = 1
$ws.Range("B18").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C18").Value2 = "Arica y Parinacota"
$ws.Range("D18").Value2 = 45054
$ws.Range("E18").Value2 = 15
$ws.Range("F18").Value2 = 100112045
$ws.Range("G18").Value2 = "Zapallo"
$ws.Range("H18").Value2 = "Camote"
$ws.Range("I18").Value2 = "3a (cosecha)"
$ws.Range("J18").Value2 = 650
$ws.Range("K18").Value2 = 280
$ws.Range("L18").Value2 = 300
$ws.Range("M18").Value2 = 289
$ws.Range("N18").Value2 = "`$/kilo (volumen en unidades)"
$ws.Range("O18").Value2 = "Región Metropolitana"
$ws.Range("P18").Value2 = 289
$ws.Range("Q18").Value2 = 1
$ws.Range("R18").Value2 = "Hortaliza"
